$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) Replace the "m:if self.name = 'anydsl'" Word field (fldChar begin /
#    instrText ... / fldChar end) with plain literal text runs that spell
#    out the same tag wrapped in curly braces: {m:if self.name = 'anydsl'}
# ---------------------------------------------------------------------------
$f = $d.Fields.Item(1)
$fieldStart = $f.Code.Start - 1
$f.Delete()

$apos = [char]39
$openXml = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
  '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml" pkg:padding="512">' +
  '<pkg:xmlData>' +
  '<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
  '<w:body><w:p>' +
  '<w:r><w:t xml:space="preserve">{m:if </w:t></w:r>' +
  '<w:r><w:t xml:space="preserve">self.name </w:t></w:r>' +
  '<w:r><w:t>=</w:t></w:r>' +
  '<w:r><w:t xml:space="preserve"> </w:t></w:r>' +
  "<w:r><w:t>$apos</w:t></w:r>" +
  '<w:r><w:t>anydsl</w:t></w:r>' +
  "<w:r><w:t>$apos}</w:t></w:r>" +
  '</w:p></w:body></w:document>' +
  '</pkg:xmlData></pkg:part></pkg:package>'

$r = $d.Range($fieldStart, $fieldStart)
$r.InsertXML($openXml)

# ---------------------------------------------------------------------------
# 2) Drop the trailing 4-space run that followed the "Unexpected tag EOF ..."
#    message at the end of the last paragraph.
# ---------------------------------------------------------------------------
$lastPara = $d.Paragraphs.Last
$paraEnd = $lastPara.Range.End
$trailing = $d.Range($paraEnd - 1 - 4, $paraEnd - 1)
$trailing.Delete()
